# Applies three small textual corrections to tc_p056r.docx, matching the
# target XML diff:
#   1. " la laisser deulx ou trois iours puys la gectes"
#        -> " la laisses deulx ou trois jours puys la gectes"
#        (the corrected "s" and "j" land in their own, colour-less runs,
#         same as the rest of the document's ad-hoc inserted characters)
#   2. "monstres obscur &"   -> "monstrer obscur &"
#        (the corrected "r" lands in its own, colour-less run)
#   3. "penetrante &"        -> "penetrant &"
#        (plain in-place text fix, no run split)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Grab a template "no local colour" run already present in the document
# (the lone "&" in "la plus claire &") so newly-inserted single
# characters can reuse its run formatting (i.e. everything except the
# explicit <w:color>). This reproduces the pattern the diff shows for
# the freshly corrected letters.
# ---------------------------------------------------------------------
$rngTemplate = $d.Content
$null = $rngTemplate.Find.Execute("la plus claire &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$templateRun = $d.Range($rngTemplate.End - 1, $rngTemplate.End)

# NOTE: $templateRun.FormattedText returns a live-linked Range, not an
# independent copy -- mutating its .Text (even via a variable) mutates
# $templateRun itself. So: copy the formatting onto the destination
# range first (this brings along the template's current text, "&"),
# and only then overwrite the destination's own .Text -- that leaves
# $templateRun/its FormattedText completely untouched.
function Set-ColorlessChar($range, [string]$char) {
    $range.FormattedText = $templateRun.FormattedText
    $range.Text = $char
}

# ---------------------------------------------------------------------
# 1. " la laisser deulx ou trois iours puys la gectes"
#    "laisser" -> "laisses" and "iours" -> "jours"
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(" la laisser deulx ou trois iours puys la gectes", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'la laisser ... iours' passage"
}

$base1 = $rng1.Start
$rPos = $base1 + 10   # the "r" in "laisser" (" la laisse|r|...")
$iPos = $base1 + 27   # the "i" in "iours"    ("...trois |i|ours...")

# Replace the later occurrence first so the earlier offset stays valid.
$iRange = $d.Range($iPos, $iPos + 1)
if ($iRange.Text -ne "i") { throw "Unexpected character at iPos: [$($iRange.Text)]" }
Set-ColorlessChar $iRange "j"

$rRange = $d.Range($rPos, $rPos + 1)
if ($rRange.Text -ne "r") { throw "Unexpected character at rPos: [$($rRange.Text)]" }
Set-ColorlessChar $rRange "s"

# ---------------------------------------------------------------------
# 2. "monstres obscur &" -> "monstrer obscur &"
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("monstres obscur &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'monstres obscur &'"
}
$sPos = $rng2.Start + 7   # the "s" in "monstres"
$sRange = $d.Range($sPos, $sPos + 1)
if ($sRange.Text -ne "s") { throw "Unexpected character at sPos: [$($sRange.Text)]" }
Set-ColorlessChar $sRange "r"

# ---------------------------------------------------------------------
# 3. "penetrante &" -> "penetrant &"  (plain text fix, same run)
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("penetrante &", $true, $false, $false, $false, $false, $true, 1, $false, `
    "penetrant &", 2)
if (-not $found3) {
    throw "Could not find/replace 'penetrante &'"
}
